# Backup before MoClo code restructure
# - Fix Transfer Volume values in rows 4 and 5
# - Add two new data rows (6 and 7) for wells A5 and A6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct existing transfer volumes
$ws.Range("H4").Value = 1875
$ws.Range("H5").Value = 2875

# Row 6: UID 5, Destination Well A5
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "level 1 6RES source plate"
$ws.Range("C6").Value = "6RES_AQ_BP"
$ws.Range("D6").Value = "A1"
$ws.Range("E6").Value = "384-Well Level 1 MoClo output plate"
$ws.Range("F6").Value = "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)"
$ws.Range("G6").Value = "A5"
$ws.Range("H6").Value = 2625
$ws.Range("I6").Value = "Deionised water"

# Row 7: UID 6, Destination Well A6
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "level 1 6RES source plate"
$ws.Range("C7").Value = "6RES_AQ_BP"
$ws.Range("D7").Value = "A1"
$ws.Range("E7").Value = "384-Well Level 1 MoClo output plate"
$ws.Range("F7").Value = "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)"
$ws.Range("G7").Value = "A6"
$ws.Range("H7").Value = 1875
$ws.Range("I7").Value = "Deionised water"
